$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (pushes the existing rows 51-52 down to 52-53)
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new weekly price record
$ws.Cells.Item(51,1).Value = 2
$ws.Cells.Item(51,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(51,3).Value = "Coquimbo"
$ws.Cells.Item(51,4).Value = 44714
$ws.Cells.Item(51,5).Value = 4
$ws.Cells.Item(51,6).Value = 100112022
$ws.Cells.Item(51,7).Value = "Arveja Verde"
$ws.Cells.Item(51,8).Value = "Perfection"
$ws.Cells.Item(51,9).Value = "Primera"
$ws.Cells.Item(51,10).Value = 240
$ws.Cells.Item(51,11).Value = 27000
$ws.Cells.Item(51,12).Value = 28000
$ws.Cells.Item(51,13).Value = 27500
$ws.Cells.Item(51,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(51,15).Value = "Provincia de Limarí"
$ws.Cells.Item(51,16).Value = 1100
$ws.Cells.Item(51,17).Value = 25
$ws.Cells.Item(51,18).Value = "Hortaliza"
